$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2033.25
$ws.Range("I62").Value = 1722.1111
$ws.Range("J62").Value = 2966.6667
$ws.Range("K62").Value = 1722.1111
$ws.Range("L62").Value = 2966.6667
$ws.Range("M62").Value = -1098.1111
$ws.Range("N62").Value = -4214.6667

$ws.Range("H64").Value = 2943.8333
$ws.Range("I64").Value = 2890
$ws.Range("J64").Value = 3051.5
$ws.Range("K64").Value = 2890
$ws.Range("L64").Value = 3051.5
$ws.Range("M64").Value = -2642
$ws.Range("N64").Value = -3547.5

$ws.Range("H65").Value = 2033.25
$ws.Range("I65").Value = 1722.1111
$ws.Range("J65").Value = 2966.6667
$ws.Range("K65").Value = 8610.5555
$ws.Range("L65").Value = 14833.3335
$ws.Range("M65").Value = -5490.5555
$ws.Range("N65").Value = -21073.3335

$ws.Range("H67").Value = 2943.8333
$ws.Range("I67").Value = 2890
$ws.Range("J67").Value = 3051.5
$ws.Range("K67").Value = 2890
$ws.Range("L67").Value = 3051.5
$ws.Range("M67").Value = -2032
$ws.Range("N67").Value = -4767.5

$ws.Range("H93").Value = 50007270
$ws.Range("I93").Value = 14546
$ws.Range("J93").Value = 100000000
$ws.Range("K93").Value = 14546
$ws.Range("L93").Value = 100000000
$ws.Range("M93").Value = -12050
$ws.Range("N93").Value = -100004992

$ws.Range("H107").Value = 536.4074000000001
$ws.Range("I107").Value = 384.29166
$ws.Range("K107").Value = 384.29166
$ws.Range("M107").Value = 1535.70834

$ws.Range("H116").Value = 2106993.5
$ws.Range("I116").Value = 11906930
$ws.Range("K116").Value = 11906930
$ws.Range("M116").Value = -11903488

$ws.Range("H123").Value = 24980
$ws.Range("J123").Value = 24980
$ws.Range("L123").Value = 24980
$ws.Range("N123").Value = -34780

$ws.Range("H130").Value = 78000
$ws.Range("J130").Value = 78000
$ws.Range("L130").Value = 78000
$ws.Range("N130").Value = -88040

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7678.9155
$ws.Range("I32").Value = 4550.0923
$ws.Range("J32").Value = 18977.445
$ws.Range("K32").Value = 4550.0923
$ws.Range("L32").Value = 18977.445
$ws.Range("M32").Value = -4263.0923
$ws.Range("N32").Value = -19551.445

$ws.Range("H45").Value = 58825492
$ws.Range("I45").Value = 71429960
$ws.Range("J45").Value = 4666.6665
$ws.Range("K45").Value = 71429960
$ws.Range("L45").Value = 4666.6665
$ws.Range("M45").Value = -71429583
$ws.Range("N45").Value = -5420.6665

$ws.Range("H92").Value = 31000
$ws.Range("J92").Value = 31000
$ws.Range("L92").Value = 31000
$ws.Range("N92").Value = -35992

$ws.Range("H102").Value = 1529.4117
$ws.Range("I102").Value = 1529.4117
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1529.4117
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 92.58829999999989
$ws.Range("N102").ClearContents()

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 2617.8276
$ws.Range("I122").Value = 2465.5
$ws.Range("K122").Value = 7396.5
$ws.Range("M122").Value = -4946.5

$ws.Range("H132").Value = 1228
$ws.Range("I132").Value = 1246.3125
$ws.Range("J132").Value = 1081.5
$ws.Range("K132").Value = 3738.9375
$ws.Range("L132").Value = 3244.5
$ws.Range("M132").Value = -1208.9375
$ws.Range("N132").Value = -8304.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1598.6957
$ws.Range("I99").Value = 1461.579
$ws.Range("J99").Value = 2250
$ws.Range("K99").Value = 1461.579
$ws.Range("L99").Value = 2250
$ws.Range("M99").Value = 36.42100000000005
$ws.Range("N99").Value = -5246

$ws.Range("H107").Value = 1088.2142
$ws.Range("I107").Value = 1020.1818
$ws.Range("J107").Value = 1337.6666
$ws.Range("K107").Value = 1020.1818
$ws.Range("L107").Value = 1337.6666
$ws.Range("M107").Value = 899.8182
$ws.Range("N107").Value = -5177.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1974.0704
$ws.Range("I31").Value = 1142.6
$ws.Range("J31").Value = 4832.25
$ws.Range("K31").Value = 1142.6
$ws.Range("L31").Value = 4832.25
$ws.Range("M31").Value = -847.5999999999999
$ws.Range("N31").Value = -5422.25

$ws.Range("H34").Value = 1974.0704
$ws.Range("I34").Value = 1142.6
$ws.Range("J34").Value = 4832.25
$ws.Range("K34").Value = 1142.6
$ws.Range("L34").Value = 4832.25
$ws.Range("M34").Value = -940.5999999999999
$ws.Range("N34").Value = -5236.25

$ws.Range("H58").Value = 835.61017
$ws.Range("I58").Value = 703.46
$ws.Range("J58").Value = 1569.7778
$ws.Range("K58").Value = 703.46
$ws.Range("L58").Value = 1569.7778
$ws.Range("M58").Value = -500.46
$ws.Range("N58").Value = -1975.7778

$ws.Range("H107").Value = 1275.8823
$ws.Range("I107").Value = 546.36365
$ws.Range("J107").Value = 2613.3333
$ws.Range("K107").Value = 546.36365
$ws.Range("L107").Value = 2613.3333
$ws.Range("M107").Value = 1373.63635
$ws.Range("N107").Value = -6453.3333

$ws.Range("H132").Value = 1302.0566
$ws.Range("I132").Value = 1081.875
$ws.Range("J132").Value = 1979.5385
$ws.Range("K132").Value = 3245.625
$ws.Range("L132").Value = 5938.6155
$ws.Range("M132").Value = -715.625
$ws.Range("N132").Value = -10998.6155

$ws.Range("H134").Value = 818.9286
$ws.Range("I134").Value = 678.3778
$ws.Range("J134").Value = 1393.909
$ws.Range("K134").Value = 2035.1334
$ws.Range("L134").Value = 4181.727000000001
$ws.Range("M134").Value = 499.8666000000001
$ws.Range("N134").Value = -9251.727000000001

$ws.Range("H136").Value = 835.61017
$ws.Range("I136").Value = 703.46
$ws.Range("J136").Value = 1569.7778
$ws.Range("K136").Value = 2110.38
$ws.Range("L136").Value = 4709.3334
$ws.Range("M136").Value = 439.6199999999999
$ws.Range("N136").Value = -9809.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 543.6
$ws.Range("J98").Value = 609
$ws.Range("L98").Value = 1827
$ws.Range("N98").Value = -4823

$ws.Range("H120").Value = 18333.334
$ws.Range("J120").Value = 18333.334
$ws.Range("L120").Value = 55000.00199999999
$ws.Range("N120").Value = -64676.00199999999

$ws.Range("H131").Value = 981.7714
$ws.Range("J131").Value = 1085.4286
$ws.Range("L131").Value = 3256.2858
$ws.Range("N131").Value = -13336.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5701.4116
$ws.Range("I70").Value = 5544.909
$ws.Range("J70").Value = 5988.3335
$ws.Range("K70").Value = 5544.909
$ws.Range("L70").Value = 5988.3335
$ws.Range("M70").Value = -5274.909
$ws.Range("N70").Value = -6528.3335

$ws.Range("H73").Value = 5701.4116
$ws.Range("I73").Value = 5544.909
$ws.Range("J73").Value = 5988.3335
$ws.Range("K73").Value = 5544.909
$ws.Range("L73").Value = 5988.3335
$ws.Range("M73").Value = -4608.909
$ws.Range("N73").Value = -7860.3335

$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 9021.214
$ws.Range("I93").Value = 26224.75
$ws.Range("J93").Value = 2139.8
$ws.Range("K93").Value = 26224.75
$ws.Range("L93").Value = 2139.8
$ws.Range("M93").Value = -24976.75
$ws.Range("N93").Value = -4635.8

$ws.Range("H122").Value = 3294.0588
$ws.Range("I122").Value = 3051.6553
$ws.Range("J122").Value = 4700
$ws.Range("K122").Value = 9154.965899999999
$ws.Range("L122").Value = 14100
$ws.Range("M122").Value = -6704.965899999999
$ws.Range("N122").Value = -19000

$ws.Range("H132").Value = 2268.52
$ws.Range("I132").Value = 1335.8235
$ws.Range("J132").Value = 4250.5
$ws.Range("K132").Value = 4007.4705
$ws.Range("L132").Value = 12751.5
$ws.Range("M132").Value = -1477.4705
$ws.Range("N132").Value = -17811.5

$ws.Range("H133").Value = 23900
$ws.Range("J133").Value = 23900
$ws.Range("L133").Value = 23900
$ws.Range("N133").Value = -28960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 13697.5
$ws.Range("J101").Value = 13697.5
$ws.Range("L101").Value = 13697.5
$ws.Range("N101").Value = -20187.5

$ws.Range("H126").Value = 1466.2727
$ws.Range("I126").Value = 933.4286
$ws.Range("J126").Value = 2398.75
$ws.Range("K126").Value = 2800.2858
$ws.Range("L126").Value = 7196.25
$ws.Range("M126").Value = -330.2857999999997
$ws.Range("N126").Value = -12136.25

$ws.Range("H132").Value = 1149.4688
$ws.Range("I132").Value = 1171.12
$ws.Range("J132").Value = 1072.1428
$ws.Range("K132").Value = 3513.36
$ws.Range("L132").Value = 3216.4284
$ws.Range("M132").Value = -983.3599999999997
$ws.Range("N132").Value = -8276.428400000001
